# Update gh-pages to output generated at 456a3b4
# Refreshes the "想去人数" (want-to-go count) and, in a couple of spots,
# the "最低票价" (lowest ticket price) columns across the four sheets of
# the workbook to match a newer scrape of the source data.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (Exhibitions) ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value  = 815
$ws.Range("F6").Value  = 699
$ws.Range("F7").Value  = 1258
$ws.Range("F9").Value  = 861
$ws.Range("F10").Value = 715
$ws.Range("F13").Value = 382
$ws.Range("F15").Value = 1019
$ws.Range("F16").Value = 11318
$ws.Range("F17").Value = 654
$ws.Range("F18").Value = 54
$ws.Range("F22").Value = 286
$ws.Range("F23").Value = 1798
$ws.Range("F27").Value = 193
$ws.Range("F29").Value = 298
$ws.Range("F30").Value = 204
$ws.Range("F31").Value = 267
$ws.Range("F32").Value = 80
$ws.Range("F36").Value = 206
$ws.Range("F37").Value = 1196
$ws.Range("F38").Value = 49

# ---- Sheet "演出" (Performances) ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("G3").Value  = 127.5
$ws.Range("F7").Value  = 150
# Row 11 is no longer sold out: want-to-go count drops by one and the
# lowest-price column switches from the text "已售罄" to a numeric 0.
$ws.Range("F11").Value = 4441
$ws.Range("G11").Value = 0
$ws.Range("F15").Value = 65
$ws.Range("F16").Value = 325
$ws.Range("F21").Value = 2

# ---- Sheet "本地生活" (Local Life) ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 836

# ---- Sheet "全部类型" (All Types) ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 836
$ws.Range("F5").Value  = 815
$ws.Range("G6").Value  = 127.5
$ws.Range("F8").Value  = 699
$ws.Range("F9").Value  = 1258
$ws.Range("F12").Value = 150
$ws.Range("F13").Value = 861
$ws.Range("F14").Value = 715
$ws.Range("F17").Value = 1019
$ws.Range("F18").Value = 11318
$ws.Range("F20").Value = 654
$ws.Range("F21").Value = 54
$ws.Range("F23").Value = 286
$ws.Range("F24").Value = 1798
$ws.Range("F27").Value = 193
$ws.Range("F32").Value = 65
$ws.Range("F33").Value = 325
$ws.Range("F34").Value = 298
$ws.Range("F36").Value = 204
$ws.Range("F37").Value = 267
$ws.Range("F38").Value = 80
$ws.Range("F45").Value = 206
$ws.Range("F46").Value = 1196
